# Update Name of Algo
# Apply updated imputed values to columns A and B for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A9"  = -22.1599
    "A18" = -22.16750000000001
    "A20" = -21.47249999999997
    "A27" = -21.8358
    "A35" = -21.89019999999999
    "A69" = -21.61799999999999
    "A76" = -19.64919999999998
    "A78" = -21.8276
    "A82" = -21.5845
    "A83" = -21.71599999999999
    "A93" = -21.3904

    "B4"  = 8.598599999999999
    "B9"  = 6.1496
    "B11" = 5.306900000000004
    "B23" = 8.6714
    "B24" = 5.249200000000001
    "B26" = 5.339000000000002
    "B34" = 9.878300000000008
    "B35" = 6.3597
    "B48" = 5.803900000000003
    "B49" = 6.3648
    "B52" = 5.419799999999993
    "B66" = 5.988
    "B67" = 5.381499999999996
    "B78" = 5.612
    "B80" = 9.694400000000003
    "B99" = 5.631999999999997
    "B104" = 10.0212
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
